$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value to D11
$ws.Range("D11").Value = "Bạn là người can đảm, cương nghị, ương ngạch."

# Add new row 17
$ws.Range("A17").Value = "Mệnh Tý Ngọ có Thiên Khốc Thiên Hư đồng cung"
$ws.Range("B17").Value = "Thiếu thời nghèo túng, trung niên khá giả, về già giàu có."

# Update selection to match final state (G20)
$ws.Range("G20").Select()
